# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets, which contain duplicated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 153
    $ws.Range("F3").Value = 29
    $ws.Range("F4").Value = 108
}
